$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date (Overview col G) / Correspond Handoff Datetime
# for the de-de sheet (col H) - both rows that reference the 14262f76... and
# b09f2bc3... files moved from 06:12:45 to 06:13:28.
$wsOverview.Range("G2").Value = "2016-08-17 06:13:28"
$wsOverview.Range("G4").Value = "2016-08-17 06:13:28"
$wsDeDe.Range("H2").Value = "2016-08-17 06:13:28"
$wsDeDe.Range("H4").Value = "2016-08-17 06:13:28"

# Priority column (E) changed from "ht" to "mt" for every row that had "ht"
# on both the zh-cn and de-de sheets.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (col H) and Correspond Handback
# DateTime (col K) refreshed for rows that reference the 14262f76... and
# b09f2bc3... files.
$wsZhCn.Range("H2").Value = "2016-08-17 06:13:23"
$wsZhCn.Range("H4").Value = "2016-08-17 06:13:23"
$wsZhCn.Range("K2").Value = "2016-08-17 06:13:40"
$wsZhCn.Range("K4").Value = "2016-08-17 06:13:40"

# de-de sheet: Correspond Handback DateTime (col K) refreshed for the same
# two rows.
$wsDeDe.Range("K2").Value = "2016-08-17 06:13:47"
$wsDeDe.Range("K4").Value = "2016-08-17 06:13:47"
